$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.06363980524021429
$ws.Range("D2").Value = 0.02128720379556626
$ws.Range("E2").Value = 0.1310726571657526
$ws.Range("F2").Value = 0.8968467584287936
$ws.Range("G2").Value = 0.7513072746196343
$ws.Range("H2").Value = 0.7954328276419602
$ws.Range("I2").Value = 0.5477375792393744
$ws.Range("K2").Value = 0.6520098420949694
$ws.Range("M2").Value = 0.2929985920987264
$ws.Range("B3").Value = 0.05604972957387133
$ws.Range("D3").Value = 0.02104086425620721
$ws.Range("E3").Value = 0.1238348222857653
$ws.Range("F3").Value = 0.8820359146490091
$ws.Range("G3").Value = 0.7366545312302861
$ws.Range("H3").Value = 0.7939190579276669
$ws.Range("I3").Value = 0.5540500731687104
$ws.Range("K3").Value = 0.5692006543793866
$ws.Range("M3").Value = 0.2619325534778056
$ws.Range("B4").Value = 0.05138430412557682
$ws.Range("D4").Value = 0.02088919826360325
$ws.Range("E4").Value = 0.1194920486920168
$ws.Range("F4").Value = 0.8736030089369677
$ws.Range("G4").Value = 0.7282704937141204
$ws.Range("H4").Value = 0.7934644900973922
$ws.Range("I4").Value = 0.5582158832843671
$ws.Range("K4").Value = 0.5183021734767124
$ws.Range("M4").Value = 0.2429539404418648
$ws.Range("B5").Value = 0.04948199249898266
$ws.Range("D5").Value = 0.02082729872750022
$ws.Range("E5").Value = 0.1177475322208963
$ws.Range("F5").Value = 0.8703321486367486
$ws.Range("G5").Value = 0.725007258705034
$ws.Range("H5").Value = 0.7933984290393994
$ws.Range("I5").Value = 0.5599862877392798
$ws.Range("K5").Value = 0.4975473291987669
$ws.Range("M5").Value = 0.235243851911747
$ws.Range("B6").Value = 0.04916605281633224
$ws.Range("D6").Value = 0.02081701493306198
$ws.Range("E6").Value = 0.1174593706744389
$ws.Range("F6").Value = 0.8697990101575925
$ws.Range("G6").Value = 0.7244746392807713
$ws.Range("H6").Value = 0.7933946508822487
$ws.Range("I6").Value = 0.560284656297398
$ws.Range("K6").Value = 0.4941002032415724
$ws.Range("M6").Value = 0.2339650292964279
$ws.Range("B7").Value = 0.05135865317001276
$ws.Range("D7").Value = 0.02088836383455828
$ws.Range("E7").Value = 0.1194684198584213
$ws.Range("F7").Value = 0.873558227191424
$ws.Range("G7").Value = 0.7282258647834823
$ws.Range("H7").Value = 0.7934631169422062
$ws.Range("I7").Value = 0.5582394649338696
$ws.Range("K7").Value = 0.5180223200890168
$ws.Range("M7").Value = 0.2428498633197833
$ws.Range("B8").Value = 0.06102390041289141
$ws.Range("D8").Value = 0.02120235695340611
$ws.Range("E8").Value = 0.1285558644041132
$ws.Range("F8").Value = 0.8916024180757205
$ws.Range("G8").Value = 0.7461273525241694
$ws.Range("H8").Value = 0.7948121703732483
$ws.Range("I8").Value = 0.5498538972040059
$ws.Range("K8").Value = 0.6234683314498852
$ws.Range("M8").Value = 0.2822668905863495
$ws.Range("B9").Value = 0.07993069191915936
$ws.Range("D9").Value = 0.02181445853193864
$ws.Range("E9").Value = 0.1471931924692953
$ws.Range("F9").Value = 0.9322620286754244
$ws.Range("G9").Value = 0.7861318903267858
$ws.Range("H9").Value = 0.8012376666745951
$ws.Range("I9").Value = 0.5357144278073882
$ws.Range("K9").Value = 0.8298309730909068
$ws.Range("M9").Value = 0.3603436245790022
$ws.Range("B10").Value = 0.09378590041299617
$ws.Range("D10").Value = 0.02226152688188066
$ws.Range("E10").Value = 0.1614038035344478
$ws.Range("F10").Value = 0.965396704244128
$ws.Range("G10").Value = 0.818566183754541
$ws.Range("H10").Value = 0.808281245028553
$ws.Range("I10").Value = 0.5267368052176664
$ws.Range("K10").Value = 0.981216111730788
$ws.Range("M10").Value = 0.4182128269320913
$ws.Range("B11").Value = 0.1000798573902273
$ws.Range("D11").Value = 0.02246425543680175
$ws.Range("E11").Value = 0.1679853757119076
$ws.Range("F11").Value = 0.9811891202228367
$ws.Range("G11").Value = 0.8339944493690723
$ws.Range("H11").Value = 0.8119940232738827
$ws.Range("I11").Value = 0.5229601828471644
$ws.Range("K11").Value = 1.050041705912804
$ws.Range("M11").Value = 0.4446559707640176
$ws.Range("B12").Value = 0.1024617972244926
$ws.Range("D12").Value = 0.02254092351947889
$ws.Range("E12").Value = 0.1704948090001395
$ws.Range("F12").Value = 0.9872734721900684
$ws.Range("G12").Value = 0.8399345308169757
$ws.Range("H12").Value = 0.8134733976994255
$ws.Range("I12").Value = 0.5215743783005991
$ws.Range("K12").Value = 1.076098571636635
$ws.Range("M12").Value = 0.4546867378685135
$ws.Range("B13").Value = 0.1019488708694354
$ws.Range("D13").Value = 0.02252441628438362
$ws.Range("E13").Value = 0.1699535908212724
$ws.Range("F13").Value = 0.9859584571835995
$ws.Range("G13").Value = 0.8386508677867255
$ws.Range("H13").Value = 0.8131515177150561
$ws.Range("I13").Value = 0.5218708629552111
$ws.Range("K13").Value = 1.070487022988004
$ws.Range("M13").Value = 0.4525256569365581
$ws.Range("B14").Value = 0.1002758509603581
$ws.Range("D14").Value = 0.02247056503113853
$ws.Range("E14").Value = 0.1681914831852822
$ws.Range("F14").Value = 0.9816875931934987
$ws.Range("G14").Value = 0.8344811798760645
$ws.Range("H14").Value = 0.8121142590823638
$ws.Range("I14").Value = 0.5228452829116073
$ws.Range("K14").Value = 1.052185538589072
$ws.Range("M14").Value = 0.4454808589223944
$ws.Range("B15").Value = 0.09925088466147258
$ws.Range("D15").Value = 0.02243756622051052
$ws.Range("E15").Value = 0.167114381673052
$ws.Range("F15").Value = 0.9790851435185459
$ws.Range("G15").Value = 0.8319398784763621
$ws.Range("H15").Value = 0.8114884789942494
$ws.Range("I15").Value = 0.5234479187977819
$ws.Range("K15").Value = 1.040974580655416
$ws.Range("M15").Value = 0.4411679825613959
$ws.Range("B16").Value = 0.09337438268481435
$ws.Range("D16").Value = 0.02224826441572603
$ws.Range("E16").Value = 0.1609760638642612
$ws.Range("F16").Value = 0.9643791641698556
$ws.Range("G16").Value = 0.817571540222815
$ws.Range("H16").Value = 0.8080488654447322
$ws.Range("I16").Value = 0.5269898115687468
$ws.Range("K16").Value = 0.9767173779602842
$ws.Range("M16").Value = 0.416487111115444
$ws.Range("B17").Value = 0.08976694964799492
$ws.Range("D17").Value = 0.02213196282342267
$ws.Range("E17").Value = 0.1572406057613449
$ws.Range("F17").Value = 0.9555422274584373
$ws.Range("G17").Value = 0.808930194748541
$ws.Range("H17").Value = 0.8060692470702691
$ws.Range("I17").Value = 0.5292414548057529
$ws.Range("K17").Value = 0.9372872752949775
$ws.Range("M17").Value = 0.4013766924944235
$ws.Range("B18").Value = 0.08769122710648958
$ws.Range("D18").Value = 0.02206500896397401
$ws.Range("E18").Value = 0.1551030726287905
$ws.Range("F18").Value = 0.9505271172802736
$ws.Range("G18").Value = 0.8040232863291124
$ws.Range("H18").Value = 0.8049784792604271
$ws.Range("I18").Value = 0.5305654615625492
$ws.Range("K18").Value = 0.9146043386093936
$ws.Range("M18").Value = 0.3926966803609417
$ws.Range("B19").Value = 0.08698828715341733
$ws.Range("D19").Value = 0.02204232943971363
$ws.Range("E19").Value = 0.1543812215016516
$ws.Range("F19").Value = 0.9488406876798052
$ws.Range("G19").Value = 0.8023727486333314
$ws.Range("H19").Value = 0.8046173746879788
$ws.Range("I19").Value = 0.5310187109326918
$ws.Range("K19").Value = 0.9069236393412439
$ws.Range("M19").Value = 0.3897596739640079
$ws.Range("B20").Value = 0.09015105343311802
$ws.Range("D20").Value = 0.02214434962262146
$ws.Range("E20").Value = 0.1576371106841989
$ws.Range("F20").Value = 0.9564759268046146
$ws.Range("G20").Value = 0.8098435165510978
$ws.Range("H20").Value = 0.8062750257977882
$ws.Range("I20").Value = 0.5289987694054439
$ws.Range("K20").Value = 0.9414850695126518
$ws.Range("M20").Value = 0.4029840710553714
$ws.Range("B21").Value = 0.1007672980514087
$ws.Range("D21").Value = 0.02248638524844182
$ws.Range("E21").Value = 0.1687085895042273
$ws.Range("F21").Value = 0.9829392186577905
$ws.Range("G21").Value = 0.8357032590904225
$ws.Range("H21").Value = 0.8124169320296915
$ws.Range("I21").Value = 0.5225578683125178
$ws.Range("K21").Value = 1.057561288886802
$ws.Range("M21").Value = 0.4475496140695867
$ws.Range("B22").Value = 0.1076971254548198
$ws.Range("D22").Value = 0.02270933415319476
$ws.Range("E22").Value = 0.1760444546150737
$ws.Range("F22").Value = 1.000841638685074
$ws.Range("G22").Value = 0.8531741492329843
$ws.Range("H22").Value = 0.8168591398680007
$ws.Range("I22").Value = 0.518606758199855
$ws.Range("K22").Value = 1.133389420198569
$ws.Range("M22").Value = 0.4767768988712078
$ws.Range("B23").Value = 0.1039993808287676
$ws.Range("D23").Value = 0.02259039878872215
$ws.Range("E23").Value = 0.1721199154574222
$ws.Range("F23").Value = 0.9912310042837049
$ws.Range("G23").Value = 0.8437971652258227
$ws.Range("H23").Value = 0.8144489817760245
$ws.Range("I23").Value = 0.5206918559686819
$ws.Range("K23").Value = 1.092921696699307
$ws.Range("M23").Value = 0.461168382768804
$ws.Range("B24").Value = 0.08997740558505996
$ws.Range("D24").Value = 0.02213874983211639
$ws.Range("E24").Value = 0.157457819587556
$ws.Range("F24").Value = 0.9560535977744706
$ws.Range("G24").Value = 0.80943041345094
$ws.Range("H24").Value = 0.8061818457812251
$ws.Range("I24").Value = 0.5291083955674232
$ws.Range("K24").Value = 0.9395872905016915
$ws.Range("M24").Value = 0.4022573529117892
$ws.Range("B25").Value = 0.07482164901557553
$ws.Range("D25").Value = 0.02164930356603634
$ws.Range("E25").Value = 0.1420618121066113
$ws.Range("F25").Value = 0.9206930273848428
$ws.Range("G25").Value = 0.7747792773158579
$ws.Range("H25").Value = 0.7990928325783102
$ws.Range("I25").Value = 0.5392922166175893
$ws.Range("K25").Value = 0.7740465307030888
$ws.Range("M25").Value = 0.3391351353165106
